$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B3 currently holds the text "Rules String Hello (Integer hour)".
# It is being changed to the text string "1" (kept as text, not a number).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1"
